# Update cryptos.xlsx price/volume figures (and two pairs of swapped rows)
# to match the latest GitHub Actions scrape.
#
# Note: several "Price" values (column D) are plain numeric-looking strings
# (e.g. "1.005"). Excel auto-converts such text to a real number when a
# Range.Value is assigned a numeric literal, which would corrupt the exact
# text formatting of the source data. Prefixing the value with a leading
# apostrophe (') forces Excel to keep/store it as text, exactly like typing
# an apostrophe-prefixed value directly into a cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.581.08'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '1.838.66'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("D4").Value = '''1.005'
$ws.Range("E4").Value = '  -0.60%  '
$ws.Range("D5").Value = '''333.52'
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("D7").Value = '''0.4605'
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("D8").Value = '''0.3843'
$ws.Range("E8").Value = '  -1.02%  '
$ws.Range("D9").Value = '''46.42'
$ws.Range("E9").Value = '  +1.38%  '
$ws.Range("D10").Value = '''0.07860'
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("D11").Value = '''0.9637'
$ws.Range("E11").Value = '  -3.63%  '
$ws.Range("D12").Value = '''21.14'
$ws.Range("E12").Value = '  -1.72%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.830.81'
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''5.863'
$ws.Range("E14").Value = '  -1.37%  '
$ws.Range("D15").Value = '''7.091'
$ws.Range("E15").Value = '  -0.88%  '
$ws.Range("D16").Value = '''1.003'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").Value = '''89.73'
$ws.Range("E17").Value = '  +1.68%  '
$ws.Range("D18").Value = '''0.06593'
$ws.Range("E18").Value = '  -1.56%  '
$ws.Range("D19").Value = '''0.00001022'
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("D20").Value = '''17.18'
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("D21").Value = '''1.002'
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("D22").Value = '27.579.51'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '''5.328'
$ws.Range("E23").Value = '  -1.27%  '
$ws.Range("D24").Value = '''10.81'
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("D25").Value = '''2.273'
$ws.Range("E25").Value = '  -1.53%  '
$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.074.80'
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''158.40'
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").Value = '''19.44'
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").Value = '''2.051'
$ws.Range("E29").Value = '  -3.41%  '
$ws.Range("D30").Value = '''5.302'
$ws.Range("E30").Value = '  -2.32%  '
$ws.Range("D31").Value = '''118.15'
$ws.Range("E31").Value = '  -2.41%  '
$ws.Range("D32").Value = '''0.09396'
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").Value = '''0.9365'
$ws.Range("E33").Value = '  -3.52%  '
$ws.Range("D34").Value = '''3.574'
$ws.Range("E34").Value = '  -1.08%  '
$ws.Range("D35").Value = '''5.234'
$ws.Range("E35").Value = '  -1.06%  '
$ws.Range("D36").Value = '''1.322'
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("D37").Value = '''0.05960'
$ws.Range("E37").Value = '  -0.61%  '
$ws.Range("D38").Value = '''0.02187'
$ws.Range("E38").Value = '  -1.86%  '
$ws.Range("D39").Value = '''8.149'
$ws.Range("E39").Value = '  -2.03%  '
$ws.Range("D40").Value = '''1.003'
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("D41").Value = '''1.148'
$ws.Range("E41").Value = '  -2.96%  '
$ws.Range("D42").Value = '''0.5771'
$ws.Range("E42").Value = '  -2.29%  '
$ws.Range("D43").Value = '''0.1831'
$ws.Range("E43").Value = '  -1.66%  '
$ws.Range("D44").Value = '''9.967'
$ws.Range("E44").Value = '  -3.80%  '
$ws.Range("D45").Value = '''1.269'
$ws.Range("E45").Value = '  +2.24%  '
$ws.Range("D46").Value = '''0.5411'
$ws.Range("E46").Value = '  -2.87%  '
$ws.Range("D47").Value = '''11.80'
$ws.Range("E47").Value = '  -2.67%  '
$ws.Range("D48").Value = '''1.906'
$ws.Range("E48").Value = '  +0.12%  '
$ws.Range("D49").Value = '''0.06831'
$ws.Range("E49").Value = '  +1.85%  '
$ws.Range("B50").Value = 'PaxosStandard'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range("D50").Value = '''1.003'
$ws.Range("E50").Value = '  -32.83%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = '''110.99'
$ws.Range("E51").Value = '  +0.11%  '
